$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.036.89"
$ws.Range("D3").Value = "1.826.00"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4617"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.51%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07253"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8614"
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07794"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.15%  "
$ws.Range("D13").Value = "1.848.62"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.336"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.533"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.41%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008697"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.36%  "
$ws.Range("D20").Value = "27.158.80"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.154"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("D24").Value = "2.084.82"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.57"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.844"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.092"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("E29").Value = "  -2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08831"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.963"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.441"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.131"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7213"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.081"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.463"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.52%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05239"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01939"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.950"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.209"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5159"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1627"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8587"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -15.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.179"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.57%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.17"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "102.66"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06262"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.618"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.15%  "
